$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember each connector's current vertical position/size (computed against
# the board's current row heights) before the row layout changes underneath it.
$shapeCount = $ws.Shapes.Count
$tops = @()
$heights = @()
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $ws.Shapes.Item($i)
    $tops += $shp.Top
    $heights += $shp.Height
}

# Insert two new rows above row 1, pushing the whole board (and the "Tablero 2"
# twin further down in the sheet) down by two rows.
$ws.Rows("1:2").Insert()

# Title in the newly freed row 1, styled with a bigger font and the same
# yellow fill used elsewhere on the board.
$ws.Range("B1").Value = "Tablero 1"
$ws.Range("B1").Font.Size = 20
$ws.Range("B1:C1").Interior.Color = 65535
$ws.Rows(1).RowHeight = 26.25

# The two inserted rows add 26.25 + 15 = 41.25pt of height above the board;
# shift every connector down by that same amount (and pin its height back to
# the original) so each one keeps anchoring the same board cells it did
# before, just two rows lower.
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Top = $tops[$i-1] + 41.25
    $shp.Height = $heights[$i-1]
}

# Restore selection to match the authored state.
$ws.Range("B3:F3").Select()
